$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.818.37"
$ws.Range("E2").Value = "  +1.09%  "

$ws.Range("D3").Value = "1.865.83"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.028"
$ws.Range("E4").Value = "  -0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.07"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.024"
$ws.Range("E6").Value = "  -0.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4384"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3839"
$ws.Range("E8").Value = "  +2.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07457"
$ws.Range("E9").Value = "  +1.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8852"
$ws.Range("E10").Value = "  +1.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.66"
$ws.Range("E11").Value = "  +1.49%  "

$ws.Range("D12").Value = "1.846.78"
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.797"
$ws.Range("E13").Value = "  +1.87%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.532"
$ws.Range("E14").Value = "  +0.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07172"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.56"
$ws.Range("E16").Value = "  +4.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.028"
$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009113"
$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.024"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.57"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("D21").Value = "27.844.61"
$ws.Range("E21").Value = "  +1.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.307"
$ws.Range("E22").Value = "  +1.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.26"
$ws.Range("E23").Value = "  -0.51%  "

$ws.Range("D24").Value = "2.091.09"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.058"
$ws.Range("E25").Value = "  +7.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.18"
$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.84"
$ws.Range("E27").Value = "  +1.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.395"
$ws.Range("E28").Value = "  +2.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.003"
$ws.Range("E29").Value = "  +3.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "121.86"
$ws.Range("E30").Value = "  +4.35%  "

$ws.Range("E31").Value = "  +0.75%  "

$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.214"
$ws.Range("E32").Value = "  +1.99%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7705"
$ws.Range("E33").Value = "  +1.84%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.053"
$ws.Range("E34").Value = "  +5.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.588"
$ws.Range("E35").Value = "  +2.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.025"
$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.148"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01992"
$ws.Range("E38").Value = "  +1.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05306"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.874"
$ws.Range("E40").Value = "  +2.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5228"
$ws.Range("E41").Value = "  +2.05%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1689"
$ws.Range("E42").Value = "  +1.52%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.950"
$ws.Range("E43").Value = "  +4.30%  "

$ws.Range("B44").Value = "PaxosStandard"
$ws.Range("C44").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.131"
$ws.Range("E44").Value = "  -24.13%  "

$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.737"
$ws.Range("E45").Value = "  +3.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.81"
$ws.Range("E46").Value = "  +3.43%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "110.51"
$ws.Range("E47").Value = "  +1.57%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.724"
$ws.Range("E48").Value = "  +1.38%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.026"
$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4739"
$ws.Range("E50").Value = "  +2.64%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06519"
$ws.Range("E51").Value = "  +1.73%  "
